# Adds a new "Axx" / "input/output errors" error-code group:
#  - Error types sheet gets a new row (Axx | input/output errors)
#  - Error codes sheet gets two new rows describing the new errors
#  - Both tables are resized to include the new rows

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Error codes"
$ws2 = $wb.Worksheets.Item(2)   # "Error types"

# --- Error types sheet: new group row (write first so the new shared
# strings "Axx" / "input/output errors" land at the lowest indices,
# matching the order they were introduced upstream) ---
$ws2.Range("A12").Value = "Axx"
$ws2.Range("B12").Value = "input/output errors"

# --- Error codes sheet: two new rows for the new group ---
$ws1.Range("A59").Value = "The file cannot be written into the selected folder. Either it does not exist or it is blocked by another program."
$ws1.Range("B59").Value = "ERRA00"
$ws1.Range("C59").Value = "yes"

# Row 60 mirrors the existing "group header" row style (row 47, the
# WARN700 entry) which carries left/right border formatting via style
# indices 1 and 2 on columns A and B.
$ws1.Range("A47:C47").Copy() | Out-Null
$ws1.Range("A60:C60").PasteSpecial(-4122) | Out-Null
$ws1.Range("B60").Value = "WARN001"
$ws1.Range("A60").Value = "The 'check report' validation of a report is not passed (validation fired from the check report in the toolbar)"
$ws1.Range("C60").Value = "yes"

# --- Resize both tables to cover the newly added rows ---
$lo1 = $ws1.ListObjects.Item(1)
$null = $lo1.Resize($ws1.Range("A1:C60"))

$lo2 = $ws2.ListObjects.Item(1)
$null = $lo2.Resize($ws2.Range("A1:B12"))

# --- Update view selections to match the edited state ---
$ws1.Activate() | Out-Null
$ws1.Range("A57").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A11").Select() | Out-Null

$ws1.Activate() | Out-Null
